$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "নাম: Most. Kaniz Fatema Isha"
